$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "goodwill " column header (note trailing space)
$ws.Range("E1").Value = "goodwill "

# Set the amount for the new column on row 2
$ws.Range("E2").Value = 5

# Update the active selection to E2
$ws.Range("E2").Select()
